# Replace the 25 multiplication-equation answer cells in the document's
# table with their new values, matching the target commit's content.
$d = $word.ActiveDocument

$d.Content.Find.Execute("287×4=1148", $true, $false, $false, $false, $false, $true, 1, $false, "710×9=6390", 2) | Out-Null
$d.Content.Find.Execute("829×5=4145", $true, $false, $false, $false, $false, $true, 1, $false, "723×2=1446", 2) | Out-Null
$d.Content.Find.Execute("425×4=1700", $true, $false, $false, $false, $false, $true, 1, $false, "381×5=1905", 2) | Out-Null
$d.Content.Find.Execute("445×8=3560", $true, $false, $false, $false, $false, $true, 1, $false, "419×7=2933", 2) | Out-Null
$d.Content.Find.Execute("264×6=1584", $true, $false, $false, $false, $false, $true, 1, $false, "799×8=6392", 2) | Out-Null
$d.Content.Find.Execute("231×7=1617", $true, $false, $false, $false, $false, $true, 1, $false, "226×8=1808", 2) | Out-Null
$d.Content.Find.Execute("172×9=1548", $true, $false, $false, $false, $false, $true, 1, $false, "922×4=3688", 2) | Out-Null
$d.Content.Find.Execute("560×6=3360", $true, $false, $false, $false, $false, $true, 1, $false, "652×5=3260", 2) | Out-Null
$d.Content.Find.Execute("726×3=2178", $true, $false, $false, $false, $false, $true, 1, $false, "720×7=5040", 2) | Out-Null
$d.Content.Find.Execute("292×3=876", $true, $false, $false, $false, $false, $true, 1, $false, "764×7=5348", 2) | Out-Null
$d.Content.Find.Execute("727×5=3635", $true, $false, $false, $false, $false, $true, 1, $false, "814×9=7326", 2) | Out-Null
$d.Content.Find.Execute("536×5=2680", $true, $false, $false, $false, $false, $true, 1, $false, "165×4=660", 2) | Out-Null
$d.Content.Find.Execute("232×8=1856", $true, $false, $false, $false, $false, $true, 1, $false, "265×4=1060", 2) | Out-Null
$d.Content.Find.Execute("751×4=3004", $true, $false, $false, $false, $false, $true, 1, $false, "231×3=693", 2) | Out-Null
$d.Content.Find.Execute("423×4=1692", $true, $false, $false, $false, $false, $true, 1, $false, "337×9=3033", 2) | Out-Null
$d.Content.Find.Execute("857×3=2571", $true, $false, $false, $false, $false, $true, 1, $false, "522×2=1044", 2) | Out-Null
$d.Content.Find.Execute("823×2=1646", $true, $false, $false, $false, $false, $true, 1, $false, "683×5=3415", 2) | Out-Null
$d.Content.Find.Execute("497×4=1988", $true, $false, $false, $false, $false, $true, 1, $false, "893×7=6251", 2) | Out-Null
$d.Content.Find.Execute("910×9=8190", $true, $false, $false, $false, $false, $true, 1, $false, "329×9=2961", 2) | Out-Null
$d.Content.Find.Execute("112×3=336", $true, $false, $false, $false, $false, $true, 1, $false, "618×9=5562", 2) | Out-Null
$d.Content.Find.Execute("182×3=546", $true, $false, $false, $false, $false, $true, 1, $false, "253×7=1771", 2) | Out-Null
$d.Content.Find.Execute("344×2=688", $true, $false, $false, $false, $false, $true, 1, $false, "779×4=3116", 2) | Out-Null
$d.Content.Find.Execute("353×5=1765", $true, $false, $false, $false, $false, $true, 1, $false, "399×9=3591", 2) | Out-Null
$d.Content.Find.Execute("387×3=1161", $true, $false, $false, $false, $false, $true, 1, $false, "376×5=1880", 2) | Out-Null
$d.Content.Find.Execute("410×3=1230", $true, $false, $false, $false, $false, $true, 1, $false, "326×7=2282", 2) | Out-Null
